$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new applicant row (row 6) to the sheet.
# Force the whole row to Text format first so that numeric-looking values
# (phone numbers, the "060" contract number, the "21-04-2025" date string)
# are stored verbatim as text instead of being auto-coerced into numbers /
# dates and losing leading zeros / "+" prefixes. ClearFormats() afterwards
# drops the temporary "@" number-format style again so the new cells end up
# with the plain default style, matching the rest of the sheet.
$ws.Range("A6:I6").NumberFormat = "@"

$ws.Range("A6").Value = "Giyosov Azizbek Ilhomovich"
$ws.Range("B6").Value = "Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik"
$ws.Range("C6").Value = "AD0993829"
$ws.Range("D6").Value = "060"
$ws.Range("E6").Value = "Qashqadaryo viloyati"
$ws.Range("F6").Value = "Koson tumani"
$ws.Range("G6").Value = "998972903393"
$ws.Range("H6").Value = "21-04-2025"
$ws.Range("I6").Value = "+998972903393"

$ws.Range("A6:I6").ClearFormats()
